$d = $word.ActiveDocument

$titles = @(
  "A Cat, A Parrot, and a Bag of Seed",
  "Socks in the Dark",
  "Predicting Fingers"
)

foreach ($title in $titles) {
  # blank separator paragraph
  $r = $d.Paragraphs.Last.Range
  $r.Collapse(0)
  $r.InsertParagraphAfter()

  # paragraph containing the title text
  $r = $d.Paragraphs.Last.Range
  $r.Collapse(0)
  $r.InsertParagraphAfter()

  $r = $d.Paragraphs.Last.Range
  $r.Collapse(0)
  $r.InsertAfter($title)
}
